# updating figures and details needed for README
# - bump the MATLAB "create bioscope files" script reference from the
#   2022/2023 version to the 2024 version
# - bump the R "Join_discreteData" script reference from v2 to v3
# - bump the R "Join_BATS_All_with_master" script reference from v2 to v3
# - leave the workbook with the "fullGrid" sheet active/selected instead
#   of "mergedSections"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("mergedSections")
$ws2 = $wb.Worksheets.Item("fullGrid")

# mergedSections sheet: update the two script-name cells
$ws1.Range("C11").Value = "Create_bioscope_files_2024_Krista.m"
$ws1.Range("C15").Value = "Join_discreteData_v3.R"

# fullGrid sheet: update the same two scripts (this sheet repeats the
# same step list in a different layout) plus the BATS join script
$ws2.Range("C8").Value  = "Join_BATS_All_with_master_v3.R"
$ws2.Range("C14").Value = "Create_bioscope_files_2024_Krista.m"

# Move the selection/active cell on each sheet and switch the active tab
# from "mergedSections" to "fullGrid"
[void]$ws1.Range("C15").Select()
[void]$ws2.Range("C14").Select()
[void]$ws2.Activate()
